$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("travels")

$newRows = @(
    @("123_Car", "Cairo", "Minya", 1),
    @("1234_Car", "Cairo", "Minya", 1),
    @("123_Car", "Cairo", "Minya", 1),
    @("123_Car", "Cairo", "Minya", 1),
    @("123_Car", "Cairo", "Minya", 1)
)

$startRow = 59
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}
